$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "order" column header in H1, copying the style of the existing
# header cell (G1) so it matches the bold/centered header formatting,
# then set its text.
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "order"

# Fill H2:H264 with a sequential order number (1-based) for each data row.
$firstRow = 2
$lastRow = 264
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 8).Value = $r - 1
}
